# Rename the "Conta D" / "Conta E" account headers to use underscores
# instead of spaces ("Conta_D" / "Conta_E"), matching the updated
# reference-data naming convention used by the new database-integrated
# reference page.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Conta_D"
$ws.Range("C1").Value = "Conta_E"

# Leave the selection where the author left it when saving.
$ws.Range("E7").Select()
